$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "money" column (C), shifting it to D.
# Excel's Insert() duplicates the left-neighbour's (B) formatting into the new column,
# which already gets us most of the way to the target look & feel.
$ws.Columns("C").Insert()

# Approximate the original author's explicit column width for the new column C.
$ws.Columns("C").ColumnWidth = 27.64

# Header row (row 4): make the new C4 header cell look like the existing (now shifted) D4
# header cell, then set its caption. Finally rename the old header text.
$ws.Range("D4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C4").Value = "Doanh số hưởng"
$ws.Range("D4").Value = "Hoa hồng"

# Body rows (5-28): numeric, left aligned, wrapped text, matching the sibling "Hoa hồng" column.
$rng = $ws.Range("C5:C28")
$rng.NumberFormat = "#,##0"
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4108
$rng.WrapText = $true

# Footer row (29): stays blank but should carry the same numeric format as the rest of column C.
$ws.Range("C29").NumberFormat = "#,##0"

# Best-effort page setup tweak (portrait orientation) referenced by the new pageSetup element.
$ws.PageSetup.Orientation = 1

# Restore/update the active selection as recorded after the edit.
$ws.Range("D20").Select()

Write-Host "done"
